$d = $word.ActiveDocument

# 1. Update the "Data:" field date: 24/10/2017 -> 25/07/2019
$found1 = $d.Content.Find.Execute("24/10/2017   ", $false, $false, $false, $false, $false, $true, 1, $false, "25/07/2019   ", 2)
if (-not $found1) {
    throw "Could not find the date text '24/10/2017   ' to replace."
}

# 2. Replace the whole clinical notes / evolution text block with the new urology note
$old = "`n# HPMA: ASSINTOMÁTICO CV`n# AP/DA:`n- DM HÁ 13 ANOS`n- HAS HÁ 13 ANOS`n- MAE = HAS, DM  PAI=FALECIDO DE TEP`n- NEGA TABAGISMO, NEGA ETILISMO`n- ATIVIDADE FISICA 2X SEMANA`n- GASTRITE`n# EM USO DE: GLIFAGE XR 2GR 2XDIA, GLIMEPIRIDA 4 2XD, VICTOZA 1,8 1XD, HUMALOG MIX 40 2XD, LOSARTAN 50 2XD, HCTZ 25 1XD, ANLODIPINA 5 1XD, AAS 100 1XD, SINVASTATINA 10 1XD, FENOFIBRATO 100 1XD`n# EXAMES:`n-- ECO TT (17/07/15): AE34  S9 PP8  VE45X30  FE61%  DISFUNÇÃO DIASTOLICA DO VE`n-- CINTILOGRAFIA MIOCARDICA (27/07/15): SEM ISQUEMIA // FE=71%`n-- LABORAT (17/07/15): UR42   CR0,84  NA137  K4,3 CT171  HDL28  LDL98  TG225   GJ177  HBGLICADA=9,8%  TGO35  TGP55  CPK264  HB14  L6000  PL204  PSA1,45  TSH3,36   URINA1=NORMAL`n# CD:`nAO ENDOCRINO`nAUMENTO ANLODIPINA`nPOLISSONOGRAFIA"
$new = "# UROLOGIA`nRETORNO COM EXAME RNM`nPSA TOTAL:4,39`nPSA LIVRE: 1,19`nHPP: CA DE IRMAO DE PROSTATA`nRNM: PI RADS II `nCD: MANTENHO MEDICAÇÃO`nPSA`n"

$found2 = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found2) {
    throw "Could not find the clinical notes block to replace."
}

Write-Host "Date updated:" $found1
Write-Host "Block replaced:" $found2
